$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header row (merged B8:C8) like the existing C4:D4 header
$ws.Range("B8").Value = "Spreadsheet SpreadsheetResult mySpr2()"
$ws.Range("B8:C8").Merge()
$ws.Range("B8:C8").Style = $ws.Range("C4:D4").Style
$ws.Range("C8").Style = $ws.Range("D4").Style

# Row 9: Steps / ]
$ws.Range("B9").Value = "Steps"
$ws.Range("C9").Value = "]"

# Row 10: Stp / (empty with border)
$ws.Range("B10").Value = "Stp"
$ws.Range("C10").Borders.LineStyle = 1
